# "final con named states"
# Adds two new worksheets ("big graph" and "exp") containing the
# epsilon-closure / subset-construction tables used to name the DFA
# states, after the existing "Sheet1" / "Sheet2".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add "big graph" sheet right after the last existing sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsBig = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsBig.Name = "big graph"

$wsBig.Range("B2").Value = "start=19"

$wsBig.Range("B3").Value = "dfa states"
$wsBig.Range("D3").Value = "nfa states"
$wsBig.Range("E3").Value = "a"
$wsBig.Range("F3").Value = "b"
$wsBig.Range("G3").Value = "c"

$wsBig.Range("B4").Value = "{20}"
$wsBig.Range("C4").Formula = '="ec("&B4&")"'
$wsBig.Range("D4").Value = "{0,2,3,4,5,6,14,19}"
$wsBig.Range("E4").Value = "{1,7,15}"
$wsBig.Range("F4").Value = "{}"
$wsBig.Range("G4").Value = "{}"

$wsBig.Range("B5").Value = "{1,7,15}"
$wsBig.Range("C5").Formula = '="ec("&B5&")"'
$wsBig.Range("D5").Value = "{5}"
$wsBig.Range("E5").Value = "{}"
$wsBig.Range("F5").Value = "{9,17}"
$wsBig.Range("G5").Value = "{}"

$wsBig.Range("B6").Formula = "=F5"
$wsBig.Range("C6").Formula = '="ec("&B6&")"'
$wsBig.Range("D6").Value = "{12,10,13}"
$wsBig.Range("E6").Value = "{}"
$wsBig.Range("F6").Value = "{11}"
$wsBig.Range("G6").Value = "{}"

$wsBig.Range("B7").Value = "{11}"
$wsBig.Range("C7").Formula = '="ec("&B7&")"'
$wsBig.Range("D7").Value = "{10,13}"
$wsBig.Range("E7").Value = "{}"
$wsBig.Range("F7").Value = "{11}"
$wsBig.Range("G7").Value = "{}"

$wsBig.Range("B10").Value = "start=20"
$wsBig.Range("D10").Value = "nfa states"
$wsBig.Range("E10").Value = "a"
$wsBig.Range("F10").Value = "b"
$wsBig.Range("G10").Value = "c"

$wsBig.Range("B11").Value = "{20}"
$wsBig.Range("C11").Formula = '="ec("&B11&")"'
$wsBig.Range("D11").Value = "{{4},{6},{14}}"
$wsBig.Range("E11").Value = "{{1},{7},{15}}"
$wsBig.Range("F11").Value = "{}"
$wsBig.Range("G11").Value = "{}"

$wsBig.Range("B12").Value = "{{1},{7},{15}}"
$wsBig.Range("C12").Formula = '="ec("&B12&")"'
$wsBig.Range("D12").Value = "{}"

# Column widths approximate the "best fit" autosize Excel applied to these
# columns (engine quantises ColumnWidth to 1/6-character increments, so the
# nearest achievable value is used).
$wsBig.Columns.Item(2).ColumnWidth = 11.0
$wsBig.Columns.Item(3).ColumnWidth = 14.666666666666666
$wsBig.Columns.Item(4).ColumnWidth = 15.666666666666666
$wsBig.Columns.Item(5).ColumnWidth = 11.0

$wsBig.Range("C28").Select()

# ---------------------------------------------------------------------
# 2. Add "exp" sheet right after "big graph"
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$wsExp.Name = "exp"

$wsExp.Range("B2").Value = "dfa states"
$wsExp.Range("C2").Value = "e-closure"
$wsExp.Range("D2").Value = "nfa states"
$wsExp.Range("E2").Value = "a"

$wsExp.Range("B3").Value = "{0}"
$wsExp.Range("C3").Formula = '="ec("&B3&")"'
$wsExp.Range("D3").Value = "{0,5}"
$wsExp.Range("E3").Value = "{1}"

$wsExp.Range("B4").Value = "{1}"
$wsExp.Range("C4").Formula = '="ec("&B4&")"'
$wsExp.Range("D4").Value = "{1,5}"
$wsExp.Range("E4").Value = "{}"

$wsExp.Range("B7").Formula = "=B2"
$wsExp.Range("C7:E7").Formula = "=C2"
$wsExp.Range("F7").Value = "b"
$wsExp.Range("G7").Value = "c"
$wsExp.Range("H7").Value = "d"

$wsExp.Range("B8").Value = "{0}"
$wsExp.Range("C8").Formula = '="ec("&B8&")"'
$wsExp.Range("D8").Value = "{0,1,2,4,7,8}"
$wsExp.Range("E8").Value = "{3,9}"
$wsExp.Range("F8").Value = "{5}"
$wsExp.Range("G8").Value = "{}"
$wsExp.Range("H8").Value = "{}"

$wsExp.Range("B9").Value = "{3,9}"
$wsExp.Range("C9").Formula = '="ec("&B9&")"'
$wsExp.Range("D9").Value = "{3,6,7,1,2,3,8,9}"
$wsExp.Range("E9").Value = "{3,9}"
$wsExp.Range("F9").Value = "{5,10}"
$wsExp.Range("G9").Value = "{}"
$wsExp.Range("H9").Value = "{}"

$wsExp.Range("B10").Value = "{5}"
$wsExp.Range("C10").Formula = '="ec("&B10&")"'
$wsExp.Range("D10").Value = "{5,6,7,8,1,2,4}"
$wsExp.Range("E10").Value = "{3,9}"
$wsExp.Range("F10").Value = "{5}"
$wsExp.Range("G10").Value = "{}"
$wsExp.Range("H10").Value = "{}"

$wsExp.Range("B11").Value = "{5,10}"
$wsExp.Range("C11").Formula = '="ec("&B11&")"'
$wsExp.Range("D11").Value = "{5,6,7,8,1,2,4,10}"
$wsExp.Range("E11").Value = "{3,9}"
$wsExp.Range("F11").Value = "{5}"
$wsExp.Range("G11").Value = "{11}"
$wsExp.Range("H11").Value = "{}"

$wsExp.Range("B12").Value = "{11}"
$wsExp.Range("C12").Formula = '="ec("&B12&")"'
$wsExp.Range("D12").Value = "{11}"
$wsExp.Range("E12").Value = "{}"
$wsExp.Range("F12").Value = "{}"
$wsExp.Range("G12").Value = "{}"
$wsExp.Range("H12").Value = "{12}"

$wsExp.Range("B13").Value = "{12}"
$wsExp.Range("C13").Formula = '="ec("&B13&")"'
$wsExp.Range("D13").Value = "{12}"
$wsExp.Range("E13").Value = "{}"
$wsExp.Range("F13").Value = "{}"
$wsExp.Range("G13").Value = "{}"
$wsExp.Range("H13").Value = "{}"

$wsExp.Range("B15").Value = "IF nfa.accept_state in nfa_states: dfa es accepting state"

$wsExp.Columns.Item(2).ColumnWidth = 13.333333333333334
$wsExp.Columns.Item(3).ColumnWidth = 13.333333333333334
$wsExp.Columns.Item(4).ColumnWidth = 14.5

$wsExp.Range("D20").Select()

# ---------------------------------------------------------------------
# 3. "Sheet2" used to be the tab shown when the file opened; now that
#    "exp" is the last-added (active) sheet, it becomes the selected
#    tab instead. Nothing else to do: Excel clears tabSelected on the
#    sheet that loses focus automatically.
# ---------------------------------------------------------------------
